# Implements "basic in dialogue input prompting" feature.
#
# Row 4 used to be a "Frog" dialogue line ("[set-talk-sfx=take_damage]Hmm").
# It is replaced by a DIALOGUE row that prompts the player for a name.
# Row 3's and row 5's dialogue text are changed to reference the name prompt
# as well (the "Tanuki" asks for / confirms the player's name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Tanuki asks the player to enter their name -----------------
$ws.Range("C3").Value = "enter your name."

# --- Row 4: rebuilt as the actual name-prompt dialogue row --------------
# Clear the old "Frog" row contents first (B4/F4/G4 no longer exist in the
# new layout), then set the new column D ("_" = no talk sfx) BEFORE column C
# (the prompt text) so that shared-string "_" is registered ahead of
# "[prompt=name]." -- matching the order new strings were appended upstream.
$ws.Range("B4").Value = ""
$ws.Range("D4").Value = "_"
$ws.Range("C4").Value = "[prompt=name]."
$ws.Range("E4").Value = "END_DIALOGUE"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""

# --- Row 5: Tanuki reacts with confirmation of the entered name ---------
$ws.Range("C5").Value = "good name."
$ws.Range("D5").Value = "_"
$ws.Range("F5").Value = "RIGHT"
$ws.Range("G5").Value = "frog_mario"
$ws.Range("H5").Value = "1, 1"
$ws.Range("I5").Value = "END_DIALOGUE"

# Update the last-selected cell, as recorded in the sheet view.
$ws.Range("C11").Select()
